$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.178.11"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.49%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.813.35"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.08%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.07"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.74%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "163.90"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.63%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.810.20"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.08%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.536"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.85%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.170"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.78%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.31"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.73%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.461"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.63%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.26"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.64%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000246"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.14%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.449.92"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.08%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.817.08"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.20%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.298.71"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.31%  "

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.56%  "

$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.55"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +5.69%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.38"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.05%  "

$ws.Range("B21").Value = "TRON"
$ws.Range("C21").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.114"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.24%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "489.50"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.69%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.29%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000158"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.46%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.91"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.87%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.26"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -4.02%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.26"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.58%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.06"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.95%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.01%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.19%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.03"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.03%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.40"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -4.79%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.959.44"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.04%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "32.02"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.34%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.760.27"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.45%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.94%  "

$ws.Range("B37").Value = "Mantle"
$ws.Range("C37").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.02"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.40%  "

$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.140"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +4.61%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.92"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.09%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.10%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.321"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.81%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.04"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.36%  "

$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "48.65"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.14%  "

$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.00"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.56%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "425.39"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -4.15%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.02%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.39"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.51%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.833.85"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.42%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "141.12"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.12%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "39.51"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.20%  "

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.55%  "
